$d = $word.ActiveDocument

# The document contains three occurrences of a split "<id>...</id>" tag
# whose middle run holds a placeholder value (p015r_a1/p015r_a2/p015r_a3).
# Replace each with the real id, collapsing the three runs into one.
$d.Content.Find.Execute("<id>p015r_a1</id>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<id>p015r_1</id>", 2)

$d.Content.Find.Execute("<id>p015r_a2</id>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<id>p015r_2</id>", 2)

$d.Content.Find.Execute("<id>p015r_a3</id>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<id>p015r_3</id>", 2)
